$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round 8 data correction: swap the "Data availability" flags for the
# Brisbane Roar Youth (row 7) and Gold Coast United (row 9) fixtures.
$ws.Range("G7").Value = "N"
$ws.Range("G9").Value = "Y"

# Move the active selection to H1, as left by the author after the edit.
$ws.Range("H1").Select()
